$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings (e.g. "4.70", "35.073.60") are not
# auto-coerced to numbers by Excel - force Text format on all target cells first.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.073.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.44"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.74"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +8.80%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0989"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.124.47"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.13"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.681"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.039.12"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.39"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.57"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.76"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.52"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.83"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +23.47%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.56%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.03"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +13.71%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +23.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.783"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +12.55%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +13.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "91.84"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.351.11"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.07"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.25%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.79"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +61.97%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0541"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.42"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.036.95"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0682"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.44"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +16.00%  "
